$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb1"
$ws.Range("C2").Value = "Ephb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.31211433333333
$ws.Range("H2").Value = 30.936343
$ws.Range("I2").Value = 0.633340936097251
$ws.Range("J2").Value = 0.633340936097251
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07588399999999999
$ws.Range("N2").Value = 0.227652
$ws.Range("O2").Value = 0.005197207581907009
$ws.Range("P2").Value = 0.00519720758190701
$ws.Range("Q2").Value = 0.7825244840706667
$ws.Range("R2").Value = 7.042720356636
$ws.Range("S2").Value = 0.003291604315016715
$ws.Range("T2").Value = 0.003291604315016716

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb1"
$ws.Range("C3").Value = "Ephb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.31211433333333
$ws.Range("H3").Value = 30.936343
$ws.Range("I3").Value = 0.633340936097251
$ws.Range("J3").Value = 0.633340936097251
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.71993
$ws.Range("N3").Value = 35.15979
$ws.Range("O3").Value = 0.80268447967186
$ws.Range("P3").Value = 0.80268447967186
$ws.Range("Q3").Value = 120.8572581386633
$ws.Range("R3").Value = 1087.71532324797
$ws.Range("S3").Value = 0.5083729397461106
$ws.Range("T3").Value = 0.5083729397461106

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb1"
$ws.Range("C4").Value = "Ephb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.31211433333333
$ws.Range("H4").Value = 30.936343
$ws.Range("I4").Value = 0.633340936097251
$ws.Range("J4").Value = 0.633340936097251
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.805103666666666
$ws.Range("N4").Value = 8.415310999999999
$ws.Range("O4").Value = 0.1921183127462331
$ws.Range("P4").Value = 0.1921183127462331
$ws.Range("Q4").Value = 28.92654972751922
$ws.Range("R4").Value = 260.338947547673
$ws.Range("S4").Value = 0.1216763920361237
$ws.Range("T4").Value = 0.1216763920361237

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb1"
$ws.Range("C5").Value = "Ephb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.103438
$ws.Range("H5").Value = 12.310314
$ws.Range("I5").Value = 0.2520215719230645
$ws.Range("J5").Value = 0.2520215719230645
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07588399999999999
$ws.Range("N5").Value = 0.227652
$ws.Range("O5").Value = 0.005197207581907009
$ws.Range("P5").Value = 0.00519720758190701
$ws.Range("Q5").Value = 0.3113852891919999
$ws.Range("R5").Value = 2.802467602728
$ws.Range("S5").Value = 0.001309808424402674
$ws.Range("T5").Value = 0.001309808424402674

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb1"
$ws.Range("C6").Value = "Ephb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.103438
$ws.Range("H6").Value = 12.310314
$ws.Range("I6").Value = 0.2520215719230645
$ws.Range("J6").Value = 0.2520215719230645
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.71993
$ws.Range("N6").Value = 35.15979
$ws.Range("O6").Value = 0.80268447967186
$ws.Range("P6").Value = 0.80268447967186
$ws.Range("Q6").Value = 48.09200611934
$ws.Range("R6").Value = 432.8280550740599
$ws.Range("S6").Value = 0.2022938043251493
$ws.Range("T6").Value = 0.2022938043251492

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb1"
$ws.Range("C7").Value = "Ephb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.103438
$ws.Range("H7").Value = 12.310314
$ws.Range("I7").Value = 0.2520215719230645
$ws.Range("J7").Value = 0.2520215719230645
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.805103666666666
$ws.Range("N7").Value = 8.415310999999999
$ws.Range("O7").Value = 0.1921183127462331
$ws.Range("P7").Value = 0.1921183127462331
$ws.Range("Q7").Value = 11.51056897973933
$ws.Range("R7").Value = 103.595120817654
$ws.Range("S7").Value = 0.04841795917351258
$ws.Range("T7").Value = 0.04841795917351258

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb1"
$ws.Range("C8").Value = "Ephb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.866538
$ws.Range("H8").Value = 5.599614
$ws.Range("I8").Value = 0.1146374919796846
$ws.Range("J8").Value = 0.1146374919796846
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07588399999999999
$ws.Range("N8").Value = 0.227652
$ws.Range("O8").Value = 0.005197207581907009
$ws.Range("P8").Value = 0.00519720758190701
$ws.Range("Q8").Value = 0.141640369592
$ws.Range("R8").Value = 1.274763326328
$ws.Range("S8").Value = 0.0005957948424876208
$ws.Range("T8").Value = 0.0005957948424876208

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb1"
$ws.Range("C9").Value = "Ephb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.866538
$ws.Range("H9").Value = 5.599614
$ws.Range("I9").Value = 0.1146374919796846
$ws.Range("J9").Value = 0.1146374919796846
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.71993
$ws.Range("N9").Value = 35.15979
$ws.Range("O9").Value = 0.80268447967186
$ws.Range("P9").Value = 0.80268447967186
$ws.Range("Q9").Value = 21.87569470234
$ws.Range("R9").Value = 196.88125232106
$ws.Range("S9").Value = 0.09201773560060014
$ws.Range("T9").Value = 0.09201773560060013

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efnb1"
$ws.Range("C10").Value = "Ephb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.866538
$ws.Range("H10").Value = 5.599614
$ws.Range("I10").Value = 0.1146374919796846
$ws.Range("J10").Value = 0.1146374919796846
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.805103666666666
$ws.Range("N10").Value = 8.415310999999999
$ws.Range("O10").Value = 0.1921183127462331
$ws.Range("P10").Value = 0.1921183127462331
$ws.Range("Q10").Value = 5.235832587772666
$ws.Range("R10").Value = 47.12249328995399
$ws.Range("S10").Value = 0.02202396153659683
$ws.Range("T10").Value = 0.02202396153659683

